$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value = [string]::Join(", ", $reversed)
        }
    }
}
